$d = $word.ActiveDocument

# The four "2018 Campaign Dates that use Perseus: Oct. 30-Nov. 8 and Nov.
# 29-Dec. 8" paragraphs (one of them preceded by a small <w:br/> run) each
# collapse down to a single, unformatted run reading
# "Campaign Dates that use Leo: April 14-23, May 14-23".
$oldNeedle = "Campaign Dates that use Perseus"
$newText = "Campaign Dates that use Leo: April 14-23, May 14-23"

$found = $true
while ($found) {
    $found = $false
    $count = $d.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $d.Paragraphs($i)
        $r = $p.Range
        $t = $r.Text
        if ($t -like "*$oldNeedle*") {
            $start = $r.Start
            $end = $r.End - 1

            # Wipe every run in the paragraph (including any leading
            # <w:br/> run) so nothing - text or formatting - survives.
            $clearRange = $d.Range($start, $end)
            $clearRange.Text = ""

            # Insert fresh, unformatted text into the now-empty paragraph;
            # this mints a brand-new run with no <w:rPr/> at all.
            $insertRange = $d.Range($start, $start)
            $insertRange.InsertAfter($newText)

            $found = $true
            break
        }
    }
}
